# Apply the authored change set to the workbook.
$wb = $excel.ActiveWorkbook
$excel.DisplayAlerts = $false

# --- Remove Sheet4 ---------------------------------------------------
# Its two cells ("Text for search product In cart" / the Mi TV title)
# are being replaced by new text that will live in Sheet2 column B.
# Deleting Sheet4 first frees up those shared-string slots so the new
# strings we write below reuse them instead of growing the string table.
$ws4 = $wb.Worksheets.Item("Sheet4")
$ws4.Delete()

# --- Sheet2: add a second column with the new copy -------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Range("B1").Value = "Scroll UpTo"
$ws2.Range("B2").Value = "LG 164 cm (65 inches) 4K Ultra HD Smart IPS LED TV 65UM7290PTD (Ceramic Black) (2020 Model)"

# Give the new column a wide, readable width (matches the authored
# ~88.4 character width as closely as the engine's pixel grid allows).
$ws2.Columns.Item(2).ColumnWidth = 87.6

# Update the view/selection so B2 is the active cell on Sheet2.
$ws2.Activate()
$ws2.Range("B2").Select()

# --- Sheet3: move the saved selection to A2 ---------------------------
$ws3 = $wb.Worksheets.Item("Sheet3")
$ws3.Activate()
$ws3.Range("A2").Select()
